$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Rename the BTec_Logo-Orange picture (header, first-page header = Headers(2))
# from image1.jpg to image2.jpg.
$hdr = $sec.Headers.Item(2)
$hdrShapeRange = $hdr.Range.InlineShapes.Item(1).Range
$hdrShapeRange.InlineShapes.Item(1).Name = "image2.jpg"

# Rename the Pearson logo picture in the primary (default) footer
# (Footers(1)) from image2.png to image1.png.
$ftr1 = $sec.Footers.Item(1)
$ftr1ShapeRange = $ftr1.Range.InlineShapes.Item(1).Range
$ftr1ShapeRange.InlineShapes.Item(1).Name = "image1.png"

# Rename the Pearson logo picture in the first-page footer
# (Footers(2)) from image2.png to image1.png.
$ftr2 = $sec.Footers.Item(2)
$ftr2ShapeRange = $ftr2.Range.InlineShapes.Item(1).Range
$ftr2ShapeRange.InlineShapes.Item(1).Name = "image1.png"
